$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the read resistance value (D2) from 685 to 1044
$ws.Range("D2").Value = 1044

# Update the temperature to test (A7) from 24.9 to 9
$ws.Range("A7").Value = 9

$wb.Save()
